$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 98) matching the existing table's layout:
# A: date (stored as text, matching the existing rows), B: weekday (text),
# C: hour (number), D: ranking (number).
$ws.Range("A98").NumberFormat = "@"
$ws.Range("A98").Value = "2025/10/13"
# Match the (unstyled) formatting of the other data rows rather than
# leaving the cell tagged with the temporary text-format style.
$ws.Range("A98").Style = $ws.Range("A97").Style

$ws.Range("B98").Value = "月"
$ws.Range("C98").Value = 8
$ws.Range("D98").Value = 201
